# Stuck Threshold Data.xlsx - update "Data" sheet values (Successor gradient run)
# and clear the stray cross-column average in C13; move the selection to D12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Column B ("Smooth") and D ("Rough") re-measurements ---------------
# Column C ("Uneven") values are unchanged by this run.

$ws.Range("B2").Value = 0.003
$ws.Range("D2").Value = 0.007

$ws.Range("B3").Value = 0.006
$ws.Range("D3").Value = 0.025

$ws.Range("B4").Value = 0.008
$ws.Range("D4").Value = 0.024

$ws.Range("B5").Value = 0.002
$ws.Range("D5").Value = 0.089

$ws.Range("B6").Value = 0.002
$ws.Range("D6").Value = 0.094

$ws.Range("B7").Value = 0.003
$ws.Range("D7").Value = 0.004

$ws.Range("B8").Value = 0.005
$ws.Range("D8").Value = 0.075

$ws.Range("B9").Value = 0.002
$ws.Range("D9").Value = 0.108

$ws.Range("B10").Value = 0.005
$ws.Range("D10").Value = 0.067

$ws.Range("B11").Value = 0.002
$ws.Range("D11").Value = 0.089

# --- Row 13 no longer carries the cross-column AVERAGE(B12:D12) --------
# (Excel kept the cell's style, just dropped the formula/value.)
$ws.Range("C13").ClearContents()

# --- Leave the selection where the author left it before saving --------
$ws.Range("D12").Select() | Out-Null
